$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 25207.25
$ws.Range("I12").Value = 50132
$ws.Range("J12").Value = 282.5
$ws.Range("K12").Value = 50132
$ws.Range("L12").Value = 282.5
$ws.Range("M12").Value = -49962
$ws.Range("N12").Value = -622.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 7412614
$ws.Range("I40").Value = 3909.9
$ws.Range("K40").Value = 3909.9
$ws.Range("M40").Value = -3734.9

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 4548.909
$ws.Range("J43").Value = 4717.125
$ws.Range("L43").Value = 4717.125
$ws.Range("N43").Value = -4855.125

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1852.6
$ws.Range("J58").Value = 4336
$ws.Range("L58").Value = 13008
$ws.Range("N58").Value = -13308

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 2354.25
$ws.Range("I96").Value = 496.75
$ws.Range("J96").Value = 4211.75
$ws.Range("K96").Value = 1490.25
$ws.Range("L96").Value = 12635.25
$ws.Range("M96").Value = -117.25
$ws.Range("N96").Value = -15381.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1364.1428
$ws.Range("I100").Value = 709.8
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 709.8
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -168.8
$ws.Range("N100").Value = -4082

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 1129.4
$ws.Range("J103").Value = 1499
$ws.Range("L103").Value = 4497
$ws.Range("N103").Value = -5669

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 93534.27
$ws.Range("I112").Value = 201599.8
$ws.Range("J112").Value = 61750.293
$ws.Range("K112").Value = 604799.3999999999
$ws.Range("L112").Value = 185250.879
$ws.Range("M112").Value = -603691.3999999999
$ws.Range("N112").Value = -187466.879

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 20456408
$ws.Range("I125").Value = 2296.4
$ws.Range("J125").Value = 37501500
$ws.Range("K125").Value = 20667.6
$ws.Range("L125").Value = 337513500
$ws.Range("M125").Value = -18207.6
$ws.Range("N125").Value = -337518420

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4567.797
$ws.Range("I138").Value = 3256.4092
$ws.Range("K138").Value = 9769.2276
$ws.Range("M138").Value = -4629.2276

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2342.7021
$ws.Range("I32").Value = 2168.7334
$ws.Range("K32").Value = 2168.7334
$ws.Range("M32").Value = -1881.7334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 71430420
$ws.Range("I61").Value = 71430420
$ws.Range("K61").Value = 71430420
$ws.Range("M61").Value = -71430208

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 20836306
$ws.Range("I74").Value = 28574070
$ws.Range("J74").Value = 3860.6924
$ws.Range("K74").Value = 28574070
$ws.Range("L74").Value = 3860.6924
$ws.Range("M74").Value = -28573196
$ws.Range("N74").Value = -5608.6924

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 20836306
$ws.Range("I77").Value = 28574070
$ws.Range("J77").Value = 3860.6924
$ws.Range("K77").Value = 142870350
$ws.Range("L77").Value = 19303.462
$ws.Range("M77").Value = -142865982
$ws.Range("N77").Value = -28039.462

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 6236.2144
$ws.Range("I122").Value = 5066.905
$ws.Range("K122").Value = 15200.715
$ws.Range("M122").Value = -12750.715

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H130").Value = 42900
$ws.Range("J130").Value = 42900
$ws.Range("L130").Value = 42900
$ws.Range("N130").Value = -52940

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 71430420
$ws.Range("I136").Value = 71430420
$ws.Range("K136").Value = 214291260
$ws.Range("M136").Value = -214288710

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 426.4375
$ws.Range("I7").Value = 321.36365
$ws.Range("J7").Value = 657.6
$ws.Range("K7").Value = 321.36365
$ws.Range("L7").Value = 657.6
$ws.Range("M7").Value = -208.36365
$ws.Range("N7").Value = -883.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 599
$ws.Range("I11").Value = 599
$ws.Range("K11").Value = 599
$ws.Range("M11").Value = -459

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 879.1429000000001
$ws.Range("I16").Value = 624.75
$ws.Range("J16").Value = 1218.3334
$ws.Range("K16").Value = 624.75
$ws.Range("L16").Value = 1218.3334
$ws.Range("M16").Value = -337.75
$ws.Range("N16").Value = -1792.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 6486.8237
$ws.Range("I22").Value = 7448.2856
$ws.Range("K22").Value = 7448.2856
$ws.Range("M22").Value = -7098.2856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 20010730
$ws.Range("I58").Value = 31265970
$ws.Range("J58").Value = 1413.5555
$ws.Range("K58").Value = 31265970
$ws.Range("L58").Value = 1413.5555
$ws.Range("M58").Value = -31265767
$ws.Range("N58").Value = -1819.5555

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 21785.643
$ws.Range("J60").Value = 24166.584
$ws.Range("L60").Value = 24166.584
$ws.Range("N60").Value = -25188.584

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H110").Value = 29998
$ws.Range("J110").Value = 29998
$ws.Range("L110").Value = 29998
$ws.Range("N110").Value = -38178

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 879.1429000000001
$ws.Range("I113").Value = 624.75
$ws.Range("J113").Value = 1218.3334
$ws.Range("K113").Value = 624.75
$ws.Range("L113").Value = 1218.3334
$ws.Range("M113").Value = 1545.25
$ws.Range("N113").Value = -5558.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 20010730
$ws.Range("I136").Value = 31265970
$ws.Range("J136").Value = 1413.5555
$ws.Range("K136").Value = 93797910
$ws.Range("L136").Value = 4240.666499999999
$ws.Range("M136").Value = -93795360
$ws.Range("N136").Value = -9340.666499999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 42.8
$ws.Range("I2").Value = 23.333334
$ws.Range("J2").Value = 72
$ws.Range("K2").Value = 140.000004
$ws.Range("L2").Value = 432
$ws.Range("M2").Value = -27.00000399999999
$ws.Range("N2").Value = -658

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 6324.3335
$ws.Range("J68").Value = 6323.25
$ws.Range("L68").Value = 18969.75
$ws.Range("N68").Value = -20591.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 6324.3335
$ws.Range("J71").Value = 6323.25
$ws.Range("L71").Value = 56909.25
$ws.Range("N71").Value = -65021.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1998.0834
$ws.Range("I140").Value = 1998.0834
$ws.Range("K140").Value = 5994.2502
$ws.Range("M140").Value = -814.2502000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3745.111
$ws.Range("I80").Value = 3636.4546
$ws.Range("J80").Value = 3915.8572
$ws.Range("K80").Value = 3636.4546
$ws.Range("L80").Value = 3915.8572
$ws.Range("M80").Value = -2638.4546
$ws.Range("N80").Value = -5911.8572

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3745.111
$ws.Range("I83").Value = 3636.4546
$ws.Range("J83").Value = 3915.8572
$ws.Range("K83").Value = 18182.273
$ws.Range("L83").Value = 19579.286
$ws.Range("M83").Value = -13190.273
$ws.Range("N83").Value = -29563.286

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2320.1035
$ws.Range("J16").Value = 3888.4614
$ws.Range("L16").Value = 3888.4614
$ws.Range("N16").Value = -4228.4614

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1284.8572
$ws.Range("I40").Value = 1284.8572
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1284.8572
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1148.8572
$ws.Range("N40").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 30000
$ws.Range("I61").Value = 30000
$ws.Range("K61").Value = 30000
$ws.Range("M61").Value = -29798

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H98").Value = 35000
$ws.Range("J98").Value = 35000
$ws.Range("L98").Value = 35000
$ws.Range("N98").Value = -40990

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 30000
$ws.Range("I113").Value = 30000
$ws.Range("K113").Value = 30000
$ws.Range("M113").Value = -27830

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2941.8125
$ws.Range("I122").Value = 2920.5386
$ws.Range("J122").Value = 3034
$ws.Range("K122").Value = 8761.6158
$ws.Range("L122").Value = 9102
$ws.Range("M122").Value = -6311.6158
$ws.Range("N122").Value = -14002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7818074.5
$ws.Range("I132").Value = 8626637
$ws.Range("K132").Value = 25879911
$ws.Range("M132").Value = -25877381

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1038.7778
$ws.Range("J81").Value = 900
$ws.Range("L81").Value = 1800
$ws.Range("N81").Value = -3922

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1038.7778
$ws.Range("J84").Value = 900
$ws.Range("L84").Value = 9000
$ws.Range("N84").Value = -19608

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2091.5557
$ws.Range("I96").Value = 1460
$ws.Range("K96").Value = 1460
$ws.Range("M96").Value = -87

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 11114937
$ws.Range("I132").Value = 14287074
$ws.Range("J132").Value = 12459.4
$ws.Range("K132").Value = 42861222
$ws.Range("L132").Value = 37378.2
$ws.Range("M132").Value = -42858692
$ws.Range("N132").Value = -42438.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 26318958
$ws.Range("I136").Value = 26318958
$ws.Range("K136").Value = 78956874
$ws.Range("M136").Value = -78954324
